# Auto-generated script to apply numeric corrections to Leviathan_Profits leve-profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2950.2666
$ws.Range("J121").Value = 3089.5715
$ws.Range("L121").Value = 9268.7145
$ws.Range("N121").Value = -12762.7145

$ws.Range("H137").Value = 2459.7441
$ws.Range("J137").Value = 3027.8667
$ws.Range("L137").Value = 9083.6001
$ws.Range("N137").Value = -14183.6001

$ws.Range("H138").Value = 2772.5
$ws.Range("J138").Value = 3135.4583
$ws.Range("L138").Value = 9406.374899999999
$ws.Range("N138").Value = -19686.3749

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 42865.56
$ws.Range("I32").Value = 24974.373
$ws.Range("K32").Value = 24974.373
$ws.Range("M32").Value = -24687.373

$ws.Range("H39").Value = 13100
$ws.Range("J39").Value = 27500
$ws.Range("L39").Value = 27500
$ws.Range("N39").Value = -28540

$ws.Range("H61").Value = 1107.7368
$ws.Range("I61").Value = 1124.4324
$ws.Range("J61").Value = 490
$ws.Range("K61").Value = 1124.4324
$ws.Range("L61").Value = 490
$ws.Range("M61").Value = -912.4323999999999
$ws.Range("N61").Value = -914

$ws.Range("H74").Value = 1464.0944
$ws.Range("I74").Value = 1347.8914
$ws.Range("K74").Value = 1347.8914
$ws.Range("M74").Value = -473.8914

$ws.Range("H77").Value = 1464.0944
$ws.Range("I77").Value = 1347.8914
$ws.Range("K77").Value = 6739.457
$ws.Range("M77").Value = -2371.457

$ws.Range("H101").Value = 28481.4
$ws.Range("J101").Value = 28481.4
$ws.Range("L101").Value = 28481.4
$ws.Range("N101").Value = -34971.4

$ws.Range("H132").Value = 2202.4138
$ws.Range("I132").Value = 1275.12
$ws.Range("K132").Value = 3825.36
$ws.Range("M132").Value = -1295.36

$ws.Range("H136").Value = 1107.7368
$ws.Range("I136").Value = 1124.4324
$ws.Range("J136").Value = 490
$ws.Range("K136").Value = 3373.2972
$ws.Range("L136").Value = 1470
$ws.Range("M136").Value = -823.2972
$ws.Range("N136").Value = -6570

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 1988.6666
$ws.Range("I29").Value = 1988.6666
$ws.Range("K29").Value = 1988.6666
$ws.Range("M29").Value = -1699.6666

$ws.Range("H80").Value = 1243.2667
$ws.Range("I80").Value = 819.1667
$ws.Range("J80").Value = 1526
$ws.Range("K80").Value = 819.1667
$ws.Range("L80").Value = 1526
$ws.Range("M80").Value = 178.8333
$ws.Range("N80").Value = -3522

$ws.Range("H83").Value = 1243.2667
$ws.Range("I83").Value = 819.1667
$ws.Range("J83").Value = 1526
$ws.Range("K83").Value = 4095.8335
$ws.Range("L83").Value = 7630
$ws.Range("M83").Value = 896.1665000000003
$ws.Range("N83").Value = -17614

$ws.Range("H134").Value = 1930.3572
$ws.Range("I134").Value = 1925
$ws.Range("K134").Value = 5775
$ws.Range("M134").Value = -3240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3455.85
$ws.Range("I31").Value = 1870.9333
$ws.Range("J31").Value = 8210.6
$ws.Range("K31").Value = 1870.9333
$ws.Range("L31").Value = 8210.6
$ws.Range("M31").Value = -1575.9333
$ws.Range("N31").Value = -8800.6

$ws.Range("H34").Value = 3455.85
$ws.Range("I34").Value = 1870.9333
$ws.Range("J34").Value = 8210.6
$ws.Range("K34").Value = 1870.9333
$ws.Range("L34").Value = 8210.6
$ws.Range("M34").Value = -1668.9333
$ws.Range("N34").Value = -8614.6

$ws.Range("H35").Value = 2200
$ws.Range("I35").Value = 2250
$ws.Range("J35").Value = 2000
$ws.Range("K35").Value = 2250
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = -1956
$ws.Range("N35").Value = -2588

$ws.Range("H86").Value = 10500
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 10500
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").Value = 10500
$ws.Range("N86").Value = -12746

$ws.Range("H89").Value = 10500
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 10500
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").Value = 52500
$ws.Range("N89").Value = -63732

$ws.Range("H122").Value = 4657
$ws.Range("I122").Value = 3980.5
$ws.Range("J122").Value = 5108
$ws.Range("K122").Value = 11941.5
$ws.Range("L122").Value = 15324
$ws.Range("M122").Value = -9491.5
$ws.Range("N122").Value = -20224

$ws.Range("H132").Value = 1759.2439
$ws.Range("I132").Value = 1716.1282
$ws.Range("K132").Value = 5148.3846
$ws.Range("M132").Value = -2618.3846

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 835.4583
$ws.Range("I5").Value = 403.9091
$ws.Range("J5").Value = 1200.6154
$ws.Range("K5").Value = 1211.7273
$ws.Range("L5").Value = 3601.8462
$ws.Range("M5").Value = -1099.7273
$ws.Range("N5").Value = -3825.8462

$ws.Range("H12").Value = 130.66667
$ws.Range("J12").Value = 215
$ws.Range("L12").Value = 645
$ws.Range("N12").Value = -991

$ws.Range("H51").Value = 899.93335
$ws.Range("I51").Value = 899.93335
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 2699.80005
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -2239.80005

$ws.Range("H92").Value = 495.66666
$ws.Range("I92").Value = 299.5
$ws.Range("J92").Value = 888
$ws.Range("K92").Value = 898.5
$ws.Range("L92").Value = 2664
$ws.Range("M92").Value = 349.5
$ws.Range("N92").Value = -5160

$ws.Range("H109").Value = 1075.9
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H135").Value = 835.4583
$ws.Range("I135").Value = 403.9091
$ws.Range("J135").Value = 1200.6154
$ws.Range("K135").Value = 3635.1819
$ws.Range("L135").Value = 10805.5386
$ws.Range("M135").Value = -1100.1819
$ws.Range("N135").Value = -15875.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 17375
$ws.Range("I22").Value = 10500
$ws.Range("J22").Value = 24250
$ws.Range("K22").Value = 10500
$ws.Range("L22").Value = 24250
$ws.Range("M22").Value = -9971
$ws.Range("N22").Value = -25308

$ws.Range("H40").Value = 22508
$ws.Range("I40").Value = 15016
$ws.Range("J40").Value = 30000
$ws.Range("K40").Value = 15016
$ws.Range("L40").Value = 30000
$ws.Range("M40").Value = -14865
$ws.Range("N40").Value = -30302

$ws.Range("H70").Value = 5754.8
$ws.Range("I70").Value = 4255.3335
$ws.Range("J70").Value = 8004
$ws.Range("K70").Value = 4255.3335
$ws.Range("L70").Value = 8004
$ws.Range("M70").Value = -3985.3335
$ws.Range("N70").Value = -8544

$ws.Range("H73").Value = 5754.8
$ws.Range("I73").Value = 4255.3335
$ws.Range("J73").Value = 8004
$ws.Range("K73").Value = 4255.3335
$ws.Range("L73").Value = 8004
$ws.Range("M73").Value = -3319.3335
$ws.Range("N73").Value = -9876

$ws.Range("H132").Value = 2920.3684
$ws.Range("I132").Value = 2581.6155
$ws.Range("J132").Value = 3654.3333
$ws.Range("K132").Value = 7744.8465
$ws.Range("L132").Value = 10962.9999
$ws.Range("M132").Value = -5214.8465
$ws.Range("N132").Value = -16022.9999

$ws.Range("H141").Value = 90390
$ws.Range("I141").Value = 90390
$ws.Range("K141").Value = 90390
$ws.Range("M141").Value = -85210

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 68500
$ws.Range("I7").Value = 88333.336
$ws.Range("K7").Value = 88333.336
$ws.Range("M7").Value = -88221.336

$ws.Range("H25").Value = 1669500
$ws.Range("I25").Value = 2503000
$ws.Range("J25").Value = 2500
$ws.Range("K25").Value = 2503000
$ws.Range("L25").Value = 2500
$ws.Range("M25").Value = -2502770
$ws.Range("N25").Value = -2960

$ws.Range("H35").Value = 10010.333
$ws.Range("I35").Value = 2515.5
$ws.Range("J35").Value = 25000
$ws.Range("K35").Value = 2515.5
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = -2179.5
$ws.Range("N35").Value = -25672

$ws.Range("H76").Value = 19429.334
$ws.Range("J76").Value = 19429.334
$ws.Range("L76").Value = 19429.334
$ws.Range("N76").Value = -20105.334

$ws.Range("H79").Value = 19429.334
$ws.Range("J79").Value = 19429.334
$ws.Range("L79").Value = 19429.334
$ws.Range("N79").Value = -21769.334

$ws.Range("H126").Value = 68500
$ws.Range("I126").Value = 88333.336
$ws.Range("K126").Value = 265000.008
$ws.Range("M126").Value = -262530.008

$ws.Range("H132").Value = 3076.0981
$ws.Range("I132").Value = 2251.3125
$ws.Range("J132").Value = 4465.2104
$ws.Range("K132").Value = 6753.9375
$ws.Range("L132").Value = 13395.6312
$ws.Range("M132").Value = -4223.9375
$ws.Range("N132").Value = -18455.6312

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3783.7666
$ws.Range("I81").Value = 3750.6296
$ws.Range("K81").Value = 7501.2592
$ws.Range("M81").Value = -6440.2592

$ws.Range("H84").Value = 3783.7666
$ws.Range("I84").Value = 3750.6296
$ws.Range("K84").Value = 37506.296
$ws.Range("M84").Value = -32202.296

$ws.Range("H96").Value = 19000
$ws.Range("I96").Value = 19000
$ws.Range("K96").Value = 19000
$ws.Range("M96").Value = -17627

$ws.Range("H122").Value = 1386.5
$ws.Range("I122").Value = 1352.8572
$ws.Range("K122").Value = 4058.5716
$ws.Range("M122").Value = -1608.5716

$ws.Range("H132").Value = 4331.737
$ws.Range("I132").Value = 4431.243
$ws.Range("J132").Value = 650
$ws.Range("K132").Value = 13293.729
$ws.Range("L132").Value = 1950
$ws.Range("M132").Value = -10763.729
$ws.Range("N132").Value = -7010

$ws.Range("H136").Value = 757.1429000000001
$ws.Range("I136").Value = 550
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 1650
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = 900
$ws.Range("N136").Value = -11100

$ws.Range("H138").Value = 100429
$ws.Range("J138").Value = 100429
$ws.Range("L138").Value = 100429
$ws.Range("N138").Value = -110709
